$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fix typo: "per capita" -> "per cap." in predictor labels (column C)
# Also fixes mismatched bracket on the "Livestock AB Consumption" row.

$ws.Range("C2").Value  = "ln(GDP [dollars per cap.])"
$ws.Range("C17").Value = "ln(GDP [dollars per cap.])"

$ws.Range("C4").Value  = "ln(Tourism - Inbound [per cap.])"

$ws.Range("C5").Value  = "ln(ProMed Mentions [per cap.])"
$ws.Range("C13").Value = "ln(ProMed Mentions [per cap.])"

$ws.Range("C6").Value  = "ln(Migrant Population [per cap.])"

$ws.Range("C9").Value  = "ln(AB Exports [dollars per cap.])"

$ws.Range("C11").Value = "ln(Publication Bias Index [per cap.])"
$ws.Range("C15").Value = "ln(Publication Bias Index [per cap.])"

$ws.Range("C12").Value = "Livestock AB Consumption (kg per cap.)"
